$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.669.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.712.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "672.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000235"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.714.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.711.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "474.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.655"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.861.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000128"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.14%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.168"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.55%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.702.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.50%  "
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0914"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.942"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000284"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.269"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
